$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Question 1": de-duplicate B22's style so it matches the same
# (fontId=0 / left+vcenter+indent1) xf used elsewhere, instead of the
# redundant one-off xf that only it referenced.
$ws1.Range("B22").Font.Bold = $false

# --- Sheet "Question 2": replace the old chi-square test formulas (B17/C17)
# with a fuller worked-out expected-frequency / chi-square table in columns G:J.
$ws2.Range("B17:C17").ClearContents()

# Mirror header row (cancer / without cancer / total) into H11:J11
$ws2.Range("H11").Value = "cancer"
$ws2.Range("I11").Value = "without cancer"
$ws2.Range("J11").Value = "total"

# Row 12 (smokers): per-cell chi-square contributions
$ws2.Range("G12").Value = "smokers"
$ws2.Range("H12").Formula = "=(C6-C12)^2/C12"
$ws2.Range("I12").Formula = "=(D6-D12)^2/D12"

# Row 13 (non_somker): per-cell chi-square contributions
$ws2.Range("G13").Value = "non_somker"
$ws2.Range("H13").Formula = "=(C7-C13)^2/C13"
$ws2.Range("I13").Formula = "=(D7-D13)^2/D13"

# Row 14 (total) label mirrored into column G
$ws2.Range("G14").Value = "total"

# Chi-square statistic
$ws2.Range("G17").Value = "chi"
$ws2.Range("H17").Formula = "=SUM(H12:I13)"

# Degrees of freedom
$ws2.Range("G18").Value = "df"
$ws2.Range("H18").Formula = "=2-1*2-1"

# p-value from chi-square distribution
$ws2.Range("G20").Value = "chisqr"
$ws2.Range("H20").Formula = "=_xlfn.CHISQ.DIST.RT(H17,1)"

# New helper column width
$ws2.Columns.Item(8).ColumnWidth = 13.43

# View state: make "Question 2" the active/selected sheet, zoomed to 150%,
# with C17 selected.
$ws2.Activate()
$ws2.Range("C17").Select()
$excel.ActiveWindow.Zoom = 150
